## Edit: refresh the cached "today" text of every auto-updating date
## footer field (slide master, all slide layouts, and the notes master)
## from "3/19/2019" to "22-Mar-19", and rename the
## "VersionedAddressBook" class label shape to "VersionedTravelBuddy".

$p = $ppt.ActivePresentation

$oldDate = "3/19/2019"
$newDate = "22-Mar-19"

# ppPlaceholderDate
$ppPlaceholderDate = 16

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDate = $false
        if ($sh.Type -eq 14) {  # msoPlaceholder
            try {
                if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDate = $true
                }
            } catch {
            }
        }
        if ($isDate -and $sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

# Every slide layout hanging off the master
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholders $layout.Shapes
}

# Notes master
$notesMaster = $p.NotesMaster
Update-DatePlaceholders $notesMaster.Shapes

## Rename the "VersionedAddressBook" shape label to "VersionedTravelBuddy"
$s = $p.Slides.Item(1)
$shapes = $s.Shapes
for ($i = 1; $i -le $shapes.Count; $i++) {
    $sh = $shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "VersionedAddressBook") {
            $sh.TextFrame.TextRange.Text = "VersionedTravelBuddy"
        }
    }
}
